# "Water Tower Effect Base" -- add effectspf (M) / effectcolor (N) data for
# the WATER tower sheet, rename the 1012001 unit (B6), and leave the WATER
# sheet/tab active+selected the way the author left it when they saved.

$wb = $excel.ActiveWorkbook

$fire  = $wb.Worksheets.Item("FIRE")
$water = $wb.Worksheets.Item("WATER")

# --- WATER sheet data edits --------------------------------------------

# Row 6 (1012001): renamed unit name string 해적,Pirate -> 선장,Captain
$water.Range("B6").Value = "선장,Captain"

# Row 6 gets an explicit effectspf/effectcolor pair (same cell formatting
# as the rest of that row: right aligned, like J6:L6)
$water.Range("M6").Value = 0.045
$water.Range("N6").Value = "(1,0,0,1)"
$water.Range("N6").HorizontalAlignment = -4152

# Rows 7-21: effectcolor only, default (unstyled) cells
$water.Range("N7").Value  = "(0.6,0.8,1,1)"
$water.Range("N8").Value  = "(0.6,0.8,1,1)"
$water.Range("N9").Value  = "(0.6,0.8,1,1)"
$water.Range("N10").Value = "(0.6,0.8,1,1)"
$water.Range("N11").Value = "(0.6,0.8,1,1)"
$water.Range("N12").Value = "(0.6,0.8,1,1)"
$water.Range("N13").Value = "(0.6,0.8,1,1)"
$water.Range("N14").Value = "(0.6,0.8,1,1)"
$water.Range("N15").Value = "(0.6,0.8,1,1)"
$water.Range("N16").Value = "(0.6,0.8,1,1)"
$water.Range("N17").Value = "(0.6,0.8,1,1)"
$water.Range("N18").Value = "(0.6,0.8,1,1)"
$water.Range("N19").Value = "(0.6,0.8,1,1)"
$water.Range("N20").Value = "(0.6,0.8,1,1)"
$water.Range("N21").Value = "(0.6,0.8,1,1)"

# Rows 22 & 23 (1012007 / 1013007): effectspf + effectcolor, right aligned
# like the rest of their row (J22:L22 / J23:L23)
$water.Range("M22").Value = 0.045
$water.Range("N22").Value = "(1,0,0,1)"
$water.Range("M23").Value = 0.045
$water.Range("N23").Value = "(1,0,0,1)"
$water.Range("M22:N23").HorizontalAlignment = -4152

# --- selection / active sheet -------------------------------------------
# Leave FIRE's selection where the author parked it, then activate WATER
# and select its last-edited cell so WATER ends up the active tab.
$fire.Activate()
$fire.Range("N5").Select()

$water.Activate()
$water.Range("N19").Select()
